$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The underlying "Clio Infra" GDP-per-Capita series was refreshed: the data
# series now runs through 2016 (previously 2010), and the historical figures
# were revised. Column E holds the series as text, same as the source file.
$values = @(
    "1806",
    "1841",
    "1878",
    "1911",
    "1924",
    "2007",
    "1980",
    "1986",
    "2040",
    "2109",
    "2115",
    "2131",
    "2232",
    "2106",
    "2150",
    "2153",
    "2174",
    "2271",
    "2469",
    "2699",
    "2778",
    "2895",
    "2906",
    "2986",
    "2684",
    "2238",
    "2064",
    "2013",
    "1969",
    "1937",
    "1945",
    "1940",
    "1852",
    "1694",
    "1629",
    "1466",
    "1479",
    "1562",
    "1701",
    "1776",
    "1776",
    "1719.09911990892",
    "1482.41662537559",
    "1439.28970453715",
    "1312.37202170923",
    "1165.81016002159",
    "1325.5953686561",
    "1323.69361439489",
    "1335.6607344147",
    "1301.5763001843",
    "1196.01487944941",
    "1219.38101904879",
    "1199.73312271067",
    "1155.70369460544",
    "1125.83211784106",
    "1107.53380416932",
    "1099.34178265069",
    "1067.96645622593",
    "1031.36220648499",
    "991.745235314568",
    "956.530484693982",
    "925",
    "967",
    "1011",
    "1060",
    "1103",
    "1113"
)

$firstRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("E$row").Value = "'" + $values[$i]
}

# The series now extends six more years (2011-2016); append the new rows.
$newYears = @(2011, 2012, 2013, 2014, 2015, 2016)
$startRow = 63
for ($i = 0; $i -lt $newYears.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = 508
    $ws.Range("B$row").Value = "Mozambique"
    $ws.Range("C$row").Value = "GDP per Capita"
    $ws.Range("D$row").Value = $newYears[$i]
}
